$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Gantt chart: add new task row (row 24) ---
# "5. UI design of Leave and Branch" task, assignee "A", 2024-11-18 -> 2024-11-19, status "Done"
$ws.Range("A24").Value = "5. UI design of Leave and Branch"
$ws.Range("B24").Value = "A"
$ws.Range("C24").Value = [DateTime]"2024-11-18"
$ws.Range("D24").Value = [DateTime]"2024-11-19"
$ws.Range("E24").Value = "Done"

# Match the status cell's look (fill + border) to the other "Status" cells
# in the table (e.g. E14) instead of the previously-unused blank style.
$ws.Range("E14").Copy()
$ws.Range("E24").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Selection moved to the newly edited cell ---
[void]$ws.Range("B24").Select()
